$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# Add a new "2509" period block for both workers, mirroring the existing
# 2505/2506/2507/2508 blocks (rows 16-23), and push the footer/signature
# block (old rows 28-29) down by two rows so the new data fits above it.
# ---------------------------------------------------------------------------

# Insert two fresh rows right after the current last data row (23), before
# the footer rows (28/29) -> footer becomes rows 30/31.
$ws.Rows.Item(24).Resize(2).Insert()

# New row 24 = worker 1 (EDILSON), period 2509 - same look as the other
# "interior" rows (use row 22 as the style/value template).
$ws.Range("B22:J22").Copy($ws.Range("B24:J24"))

# New row 25 = worker 2 (WILSON), period 2509 - this is now the new last
# row of the table, so it should carry the special bottom-border styling
# that used to belong to row 23.
$ws.Range("B23:J23").Copy($ws.Range("B25:J25"))

# Row 23 is no longer the last row of the table, so give it the regular
# "interior" styling (like row 22), keeping worker 2's own data.
$ws.Range("B22:J22").Copy($ws.Range("B23:J23"))
$ws.Range("C23").Value = $ws.Range("C21").Value2
$ws.Range("D23").Value = $ws.Range("D21").Value2
$ws.Range("F23").Value = $ws.Range("F21").Value2
$ws.Range("G23").Value = $ws.Range("G21").Value2

# Fill in the new period value for the two freshly-added rows.
$ws.Range("E24").Value = "2509"
$ws.Range("E25").Value = "2509"

# ---------------------------------------------------------------------------
# Update the summary figures now that a 5th period has been added.
# ---------------------------------------------------------------------------
$ws.Range("F13").Value = 5
$ws.Range("E11").Value = 609400
